# Apply updated "想去人数" (interest count) values to the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row (by sheet row number) -> new value for column F
$updates = @{
    2  = 8390
    3  = 7948
    8  = 134
    13 = 138
    14 = 1965
    19 = 131
    20 = 18
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
